$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Bảng cân đối kế toán riêng"
$ws.Range("B6").Value = "11-12"
$ws.Range("A8").Value = "Phụ lục 1 - Danh sách công ty con tại ngày 31 tháng 12 năm 2022"

# B8 must stay text ("80.83"), not be auto-coerced to a number: force text
# format before writing, then restore the cell's default style afterward.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "80.83"
$ws.Range("B8").Style = "Normal"
